$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8087598739804776
$ws.Range("C2").Value = 0.2372124282664458
$ws.Range("E2").Value = 0.1137671238965225
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.3851836419636641
$ws.Range("H2").Value = 0.5131697835815316
$ws.Range("I2").Value = 0.3537939942671571
$ws.Range("M2").Value = 0.3489902609068523
$ws.Range("B3").Value = 0.7078462491325013
$ws.Range("C3").Value = 0.2075993236583997
$ws.Range("E3").Value = 0.1087946720881305
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.3802018693666724
$ws.Range("H3").Value = 0.5171572906993589
$ws.Range("I3").Value = 0.3606072632121169
$ws.Range("M3").Value = 0.3103254445750281
$ws.Range("B4").Value = 0.645779354924116
$ws.Range("C4").Value = 0.1893471889719933
$ws.Range("E4").Value = 0.1058507841661722
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.377736924566122
$ws.Range("H4").Value = 0.5200936009997719
$ws.Range("I4").Value = 0.3652370053568212
$ws.Range("M4").Value = 0.2866612629689413
$ws.Range("B5").Value = 0.6204607587437749
$ws.Range("C5").Value = 0.1818917731582701
$ws.Range("E5").Value = 0.1046782845674663
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.3768803600719508
$ws.Range("H5").Value = 0.5214123348922328
$ws.Range("I5").Value = 0.3672352670966674
$ws.Range("M5").Value = 0.2770367868175683
$ws.Range("B6").Value = 0.6162550860845215
$ws.Range("C6").Value = 0.1806527480325713
$ws.Range("E6").Value = 0.1044852233259554
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.3767470207647676
$ws.Range("H6").Value = 0.5216386722626538
$ws.Range("I6").Value = 0.3675737988605654
$ws.Range("M6").Value = 0.2754397814913716
$ws.Range("B7").Value = 0.6454380031282767
$ws.Range("C7").Value = 0.1892467136071616
$ws.Range("E7").Value = 0.1058348618640466
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.3777247755909912
$ws.Range("H7").Value = 0.5201108919893898
$ws.Range("I7").Value = 0.3652635035819216
$ws.Range("M7").Value = 0.2865313878799753
$ws.Range("B8").Value = 0.7739871399975868
$ws.Range("C8").Value = 0.2270162407345708
$ws.Range("E8").Value = 0.1120297687232892
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.3833418114083997
$ws.Range("H8").Value = 0.51444307372374
$ws.Range("I8").Value = 0.356050182205248
$ws.Range("M8").Value = 0.3356426375975659
$ws.Range("B9").Value = 1.025216988216869
$ws.Range("C9").Value = 0.3005353470534828
$ws.Range("E9").Value = 0.1250590703536574
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.3991339319523632
$ws.Range("H9").Value = 0.5072249542675564
$ws.Range("I9").Value = 0.3415532060472124
$ws.Range("M9").Value = 0.4325721184494
$ws.Range("B10").Value = 1.209271102391028
$ws.Range("C10").Value = 0.3542292099792803
$ws.Range("E10").Value = 0.1351900343990593
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.4137400523715513
$ws.Range("H10").Value = 0.5043318589107599
$ws.Range("I10").Value = 0.3331181246827448
$ws.Range("M10").Value = 0.504198306299358
$ws.Range("B11").Value = 1.292889101552419
$ws.Range("C11").Value = 0.3785894934084695
$ws.Range("E11").Value = 0.1399247686883243
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.4210567419546294
$ws.Range("H11").Value = 0.5035465625803823
$ws.Range("I11").Value = 0.3297703958046156
$ws.Range("M11").Value = 0.5368800103832569
$ws.Range("B12").Value = 1.324537051184052
$ws.Range("C12").Value = 0.3878048544087278
$ws.Range("E12").Value = 0.1417361826034949
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.423925629884323
$ws.Range("H12").Value = 0.5033261258168551
$ws.Range("I12").Value = 0.3285737766008197
$ws.Range("M12").Value = 0.5492703687344118
$ws.Range("B13").Value = 1.317721837639169
$ws.Range("C13").Value = 0.385820576588003
$ws.Range("E13").Value = 0.1413452352023228
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.4233033707528477
$ws.Range("H13").Value = 0.503370169221725
$ws.Range("I13").Value = 0.3288283166933894
$ws.Range("M13").Value = 0.546601232858734
$ws.Range("B14").Value = 1.295493129323688
$ws.Range("C14").Value = 0.3793478338083673
$ws.Range("E14").Value = 0.1400734226291718
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.4212907898465375
$ws.Range("H14").Value = 0.5035268816549063
$ws.Range("I14").Value = 0.3296705203191017
$ws.Range("M14").Value = 0.537899079974693
$ws.Range("B15").Value = 1.281875251858992
$ws.Range("C15").Value = 0.3753818757820113
$ws.Range("E15").Value = 0.139296816506878
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.4200708624119329
$ws.Range("H15").Value = 0.5036329103173642
$ws.Range("I15").Value = 0.3301956739764158
$ws.Range("M15").Value = 0.532570659142479
$ws.Range("B16").Value = 1.203804214801096
$ws.Range("C16").Value = 0.3526358972694084
$ws.Range("E16").Value = 0.1348831748763715
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.4132755494388363
$ws.Range("H16").Value = 0.5043939128754431
$ws.Range("I16").Value = 0.3333468161973911
$ws.Range("M16").Value = 0.5020644918913035
$ws.Range("B17").Value = 1.155881718572687
$ws.Range("C17").Value = 0.3386652741210128
$ws.Range("E17").Value = 0.1322080880331455
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 0.4092801309806902
$ws.Range("H17").Value = 0.5049971448784731
$ws.Range("I17").Value = 0.3354058021268465
$ws.Range("M17").Value = 0.4833754290709606
$ws.Range("B18").Value = 1.12830773230354
$ws.Range("C18").Value = 0.3306235895069563
$ws.Range("E18").Value = 0.1306812982994501
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.407045258939803
$ws.Range("H18").Value = 0.5053940315070236
$ws.Range("I18").Value = 0.3366361012875849
$ws.Range("M18").Value = 0.4726352341110669
$ws.Range("B19").Value = 1.118969917720847
$ws.Range("C19").Value = 0.3278997554952525
$ws.Range("E19").Value = 0.1301663784133638
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 0.4062993751688992
$ws.Range("H19").Value = 0.505536966319184
$ws.Range("I19").Value = 0.3370605425551787
$ws.Range("M19").Value = 0.4690003691381008
$ws.Range("B20").Value = 1.160984214256359
$ws.Range("C20").Value = 0.3401531071751833
$ws.Range("E20").Value = 0.1324916267741827
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.4096988993376982
$ws.Range("H20").Value = 0.5049277585002159
$ws.Range("I20").Value = 0.3351818510070466
$ws.Range("M20").Value = 0.4853639510116636
$ws.Range("B21").Value = 1.302022689390242
$ws.Range("C21").Value = 0.3812492886266909
$ws.Range("E21").Value = 0.1404464807521535
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 0.4218792557769717
$ws.Range("H21").Value = 0.5034787584571205
$ws.Range("I21").Value = 0.3294212095588982
$ws.Range("M21").Value = 0.5404547168535032
$ws.Range("B22").Value = 1.394103761667907
$ws.Range("C22").Value = 0.4080535473758005
$ws.Range("E22").Value = 0.14575327043719
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.4304130210032326
$ws.Range("H22").Value = 0.5029804916732559
$ws.Range("I22").Value = 0.3260710145700898
$ws.Range("M22").Value = 0.5765444304992258
$ws.Range("B23").Value = 1.344967336497007
$ws.Range("C23").Value = 0.3937525778320605
$ws.Range("E23").Value = 0.1429109575836662
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.4258054319502236
$ws.Range("H23").Value = 0.5032051638058874
$ws.Range("I23").Value = 0.327820893850042
$ws.Range("M23").Value = 0.5572748085850492
$ws.Range("B24").Value = 1.158677446632964
$ws.Range("C24").Value = 0.3394804883452593
$ws.Range("E24").Value = 0.1323634041955586
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.4095093807009107
$ws.Range("H24").Value = 0.504958972131007
$ws.Range("I24").Value = 0.3352829543300366
$ws.Range("M24").Value = 0.4844649265785392
$ws.Range("B25").Value = 0.9573446041297871
$ws.Range("C25").Value = 0.2807037316394485
$ws.Range("E25").Value = 0.1214377892622949
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.3943408841237925
$ws.Range("H25").Value = 0.5087569279820201
$ws.Range("I25").Value = 0.3450887399219198
$ws.Range("M25").Value = 0.4062801490012333
